$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 30   Number  33"
$ws.Range("C9").Value = "Report Covering the Week  8/14/2023  Through  8/20/2023"

# --- Cells changing between numeric and placeholder-text ("0" / "***.*") ---
$ws.Range("F14").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("D15").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("C26").Value = 1
$ws.Range("I14").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("D26").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("F26").Value = 1
$ws.Range("I14").Copy()
$ws.Range("F26").PasteSpecial(-4122)
$ws.Range("G28").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("G28").PasteSpecial(-4122)
$ws.Range("H28").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("H28").PasteSpecial(-4122)
$ws.Range("G29").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("H29").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("H29").PasteSpecial(-4122)

# --- Plain numeric value updates (style/type unchanged) ---
$ws.Range("N14").Value = -61.538461538461
$ws.Range("L15").Value = -50
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = 7.692307692307
$ws.Range("I16").Value = 119
$ws.Range("J16").Value = 109
$ws.Range("K16").Value = 9.174311926605
$ws.Range("L16").Value = 46.913580246913
$ws.Range("M16").Value = -47.577092511013
$ws.Range("N16").Value = -85.416666666666
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -28.571428571428
$ws.Range("G17").Value = 30
$ws.Range("H17").Value = -33.333333333333
$ws.Range("I17").Value = 197
$ws.Range("J17").Value = 178
$ws.Range("K17").Value = 10.674157303370
$ws.Range("L17").Value = 48.120300751879
$ws.Range("M17").Value = 50.381679389313
$ws.Range("N17").Value = -44.507042253521
$ws.Range("C18").Value = 10
$ws.Range("D18").Value = 5
$ws.Range("F18").Value = 38
$ws.Range("H18").Value = 58.333333333333
$ws.Range("I18").Value = 168
$ws.Range("J18").Value = 196
$ws.Range("K18").Value = -14.285714285714
$ws.Range("L18").Value = 32.283464566929
$ws.Range("M18").Value = -40.213523131672
$ws.Range("N18").Value = -80.281690140845
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = -5.555555555555
$ws.Range("F19").Value = 70
$ws.Range("G19").Value = 74
$ws.Range("H19").Value = -5.405405405405
$ws.Range("I19").Value = 487
$ws.Range("J19").Value = 424
$ws.Range("K19").Value = 14.858490566037
$ws.Range("L19").Value = 77.737226277372
$ws.Range("M19").Value = 61.794019933554
$ws.Range("N19").Value = 46.246246246246
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 19
$ws.Range("H20").Value = -31.578947368421
$ws.Range("I20").Value = 104
$ws.Range("J20").Value = 104
$ws.Range("L20").Value = 7.216494845360
$ws.Range("M20").Value = 10.638297872340
$ws.Range("N20").Value = -81.328545780969
$ws.Range("C21").Value = 38
$ws.Range("D21").Value = 35
$ws.Range("E21").Value = 8.571428571428
$ws.Range("F21").Value = 155
$ws.Range("G21").Value = 162
$ws.Range("H21").Value = -4.320987654320
$ws.Range("I21").Value = 1087
$ws.Range("J21").Value = 1023
$ws.Range("K21").Value = 6.256109481915
$ws.Range("L21").Value = 49.108367626886
$ws.Range("M21").Value = 4.619826756496
$ws.Range("N21").Value = -63.202437373053
$ws.Range("G22").Value = 2
$ws.Range("L22").Value = 23.076923076923
$ws.Range("M22").Value = -44.827586206896
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = -25
$ws.Range("F23").Value = 15
$ws.Range("G23").Value = 18
$ws.Range("H23").Value = -16.666666666666
$ws.Range("I23").Value = 132
$ws.Range("J23").Value = 101
$ws.Range("K23").Value = 30.693069306930
$ws.Range("L23").Value = 23.364485981308
$ws.Range("M23").Value = 50
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = -13.043478260869
$ws.Range("F24").Value = 101
$ws.Range("H24").Value = 10.989010989011
$ws.Range("I24").Value = 687
$ws.Range("J24").Value = 782
$ws.Range("K24").Value = -12.148337595907
$ws.Range("L24").Value = 12.254901960784
$ws.Range("M24").Value = -12.035851472471
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 83.333333333333
$ws.Range("F25").Value = 56
$ws.Range("G25").Value = 32
$ws.Range("H25").Value = 75
$ws.Range("I25").Value = 323
$ws.Range("J25").Value = 312
$ws.Range("K25").Value = 3.525641025641
$ws.Range("L25").Value = 43.555555555555
$ws.Range("M25").Value = 3.858520900321
$ws.Range("H26").Value = -50
$ws.Range("I26").Value = 14
$ws.Range("K26").Value = -12.5
$ws.Range("L26").Value = -26.315789473684
$ws.Range("C27").Value = 3
$ws.Range("E27").Value = 200
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 40
$ws.Range("J27").Value = 30
$ws.Range("K27").Value = 33.333333333333
$ws.Range("L27").Value = -24.528301886792
$ws.Range("N28").Value = -84.615384615384
$ws.Range("N29").Value = -89.130434782608
$ws.Range("F30").Value = 1
$ws.Range("L30").Value = 0
